# Auto-generated Excel COM-interop script
# Updates the "responseAlreadyUsed" sheet with new per-category response
# texts/keys, and adds a new "priority" sheet used to rank which category
# wins when multiple intents match (used by the AI response selection).

$wb = $excel.ActiveWorkbook

# --- 1. Update responseAlreadyUsed (sheet3): columns B (key) and C (text) ---
$wsUsed = $wb.Worksheets.Item("responseAlreadyUsed")

$wsUsed.Range("B1").Value = "abuseResponseAlreadyUsed"
$wsUsed.Range("C1").Value = "[`"Just a quick reminder -- if there's any kind of abuse situation going on, or if you or anyone else is unsafe, please remember that I can't get help for you, so you would need to get help for yourself if you need it.`", `"And I'm programmed to err on the side of caution (i.e. I might accidentally give this response when it's not a clearly abusive situation, sorry about that). I didn't mean to interrupt the flow of this conversation, so please feel free to keep talking.`"]"
$wsUsed.Range("B2").Value = "imAddictedResponseAlreadyUsed"
$wsUsed.Range("C2").Value = "Addictions can be really tough. Could you say more about what it means for you?"
$wsUsed.Range("B3").Value = "imAnxiousResponseAlreadyUsed"
$wsUsed.Range("C3").Value = "If there's anything more about stress/anxiety/fear you want to explore, feel free to say more about that here"
$wsUsed.Range("B4").Value = "thisBotIsBadResponseAlreadyUsed"
$wsUsed.Range("C4").Value = "I'm a simple bot, trying my best to listen and help. I won't always get it right (sorry about that) but I'm still here to listen"
$wsUsed.Range("B5").Value = "iHateCoronavirusResponseAlreadyUsed"
$wsUsed.Range("C5").Value = "I don't know much about coronavirus or pandemics, but I can be here for you while you talk about it if that might help"
$wsUsed.Range("B6").Value = "iHaveDepressionResponseAlreadyUsed"
$wsUsed.Range("C6").Value = "I think it's sad whenever anyone feels depression or anything like that"
$wsUsed.Range("B7").Value = "feelEmptyResponseAlreadyUsed"
$wsUsed.Range("C7").Value = "I just wanted to pick up on the idea of emptiness. Feel free to say anything more about that?"
$wsUsed.Range("B8").Value = "familyProblemsResponseAlreadyUsed"
$wsUsed.Range("C8").Value = "Do keep telling me more. I'm definitely no expert on family or relationships or anything like that, but I hope talking about it is helping."
$wsUsed.Range("B9").Value = "iHateMyselfResponseAlreadyUsed"
$wsUsed.Range("C9").Value = "I'm sensing some negativity towards yourself, and that's really sad"
$wsUsed.Range("B10").Value = "helpResponseAlreadyUsed"
$wsUsed.Range("C10").Value = "If you're needing help, I can try to be useful by being a place for you to talk through what's on your mind"
$wsUsed.Range("B11").Value = "iDontKnowWhatToDoResponseAlreadyUsed"
$wsUsed.Range("C11").Value = "Hmm, sounds tough. Would you like to explore your options with me? I'll be here to listen"
$wsUsed.Range("B12").Value = "iHateHowILookResponseAlreadyUsed"
$wsUsed.Range("C12").Value = "Sometimes issues like image and appearance can make us feel bad or cause self-esteem issues. That's really sad."
$wsUsed.Range("B13").Value = "feelingLonelyResponseAlreadyUsed"
$wsUsed.Range("C13").Value = "It's sad whenever anyone is lonely or is missing the connections to other people that are so important"
$wsUsed.Range("B14").Value = "feelLostResponseAlreadyUsed"
$wsUsed.Range("C14").Value = "I'm getting the sense of a lost, almost forlorn feeling from what you're saying"
$wsUsed.Range("B15").Value = "feelOverwhelmedResponseAlreadyUsed"
$wsUsed.Range("C15").Value = "I'm sensing a certain level of feeling a bit overwhelmed, perhaps?"
$wsUsed.Range("B16").Value = "makesMeWantToSelfHarmResponseAlreadyUsed"
$wsUsed.Range("C16").Value = "I'm sorry that things have got so had that you feel that way"
$wsUsed.Range("B17").Value = "imFeelingSuicidalResponseAlreadyUsed"
$wsUsed.Range("C17").Value = "Whenever someone has suicidal thoughts, that's always sad."
$wsUsed.Range("B18").Value = "imUpsetResponseAlreadyUsed"
$wsUsed.Range("C18").Value = "I'm sorry to hear about any sadness or upset that you might be having at the moment"
$wsUsed.Range("B19").Value = "imUselessResponseAlreadyUsed"
$wsUsed.Range("C19").Value = "I hope you don't me mentioning that everyone has value, everyone has worth."

# --- 2. Add the new "priority" sheet (after responseAlreadyUsed) ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsPriority = $wb.Worksheets.Add($null, $lastSheet)
$wsPriority.Name = "priority"
$wsPriority.Columns.Item(1).ColumnWidth = 30.1

$wsPriority.Range("A1").Value = "Abuse"
$wsPriority.Range("B1").Value = 3
$wsPriority.Range("A2").Value = "Addiction"
$wsPriority.Range("B2").Value = 5
$wsPriority.Range("A3").Value = "Anxiety"
$wsPriority.Range("B3").Value = 6
$wsPriority.Range("A4").Value = "Complaining about using the bot"
$wsPriority.Range("B4").Value = 18
$wsPriority.Range("A5").Value = "Coronavirus/Lockdown"
$wsPriority.Range("B5").Value = 17
$wsPriority.Range("A6").Value = "Depression"
$wsPriority.Range("B6").Value = 2
$wsPriority.Range("A7").Value = "Empty"
$wsPriority.Range("B7").Value = 8
$wsPriority.Range("A8").Value = "Family & Relationships"
$wsPriority.Range("B8").Value = 13
$wsPriority.Range("A9").Value = "Hate myself"
$wsPriority.Range("B9").Value = 4
$wsPriority.Range("A10").Value = "Help"
$wsPriority.Range("B10").Value = 7
$wsPriority.Range("A11").Value = "I don't know what to do"
$wsPriority.Range("B11").Value = 15
$wsPriority.Range("A12").Value = "I feel ugly"
$wsPriority.Range("B12").Value = 16
$wsPriority.Range("A13").Value = "Lonely"
$wsPriority.Range("B13").Value = 11
$wsPriority.Range("A14").Value = "Lost"
$wsPriority.Range("B14").Value = 12
$wsPriority.Range("A15").Value = "Overwhelmed"
$wsPriority.Range("B15").Value = 14
$wsPriority.Range("A16").Value = "Self-harm"
$wsPriority.Range("B16").Value = 1
$wsPriority.Range("A17").Value = "Suicidal"
$wsPriority.Range("B17").Value = 0
$wsPriority.Range("A18").Value = "Upset"
$wsPriority.Range("B18").Value = 10
$wsPriority.Range("A19").Value = "Useless/Worthless/Failure"
$wsPriority.Range("B19").Value = 9

# --- 3. Restore cursor/selection state on the other sheets ---
$wsExemplars = $wb.Worksheets.Item("exemplars")
$wsExemplars.Activate() | Out-Null
$wsExemplars.Range("F2").Select() | Out-Null

$wsPriority.Activate() | Out-Null
$wsPriority.Range("G7").Select() | Out-Null

$wsUsed.Activate() | Out-Null
$wsUsed.Range("G24").Select() | Out-Null

